$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.597.95'
$ws.Range('E2').Value = '  +1.09%  '
$ws.Range('D3').Value = '2.269.03'
$ws.Range('E3').Value = '  -0.47%  '
$ws.Range('E4').Value = '  +0.33%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '120.38'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.74%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '269.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('E7').Value = '  +3.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.01'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.620'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.50'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('E11').Value = '  +1.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.39'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.39%  '
$ws.Range('E13').Value = '  -1.86%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.908'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.51%  '
$ws.Range('D16').Value = '2.612.82'
$ws.Range('E16').Value = '  -0.34%  '
$ws.Range('D17').Value = '2.270.98'
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '43.612.48'
$ws.Range('E18').Value = '  +1.03%  '
$ws.Range('E19').Value = '  +1.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.88'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.62'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.62%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.41'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '234.90'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.97'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.11%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.74'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('E26').Value = '  +8.70%  '
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.28'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.52%  '
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  +0.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '175.09'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.95%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.51'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0921'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.74%  '
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.53'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +15.84%  '
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0382'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.72'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.14%  '
$ws.Range('E39').Value = '  +4.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.56'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.96%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '13.84'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.243'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.03%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '72.45'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.73%  '
$ws.Range('E44').Value = '  +0.24%  '
$ws.Range('E45').Value = '  -1.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.75'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -7.82%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '76.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +37.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.664'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +18.23%  '
$ws.Range('E49').Value = '  +1.80%  '
$ws.Range('E50').Value = '  -1.07%  '
